$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Cells.Item(1, 1).Value = "onyen"
$ws.Cells.Item(1, 2).Value = "csid"
$ws.Cells.Item(1, 3).Value = "email"
$ws.Cells.Item(1, 4).Value = "firstName"
$ws.Cells.Item(1, 5).Value = "lastName"
$ws.Cells.Item(1, 6).Value = "pronouns"
$ws.Cells.Item(1, 7).Value = "pid"
$ws.Cells.Item(1, 8).Value = "status"
$ws.Cells.Item(1, 9).Value = "alternativeName"
$ws.Cells.Item(1, 10).Value = "gender"
$ws.Cells.Item(1, 11).Value = "ethnicity"
$ws.Cells.Item(1, 12).Value = "residency"
$ws.Cells.Item(1, 13).Value = "enteringStatus"
$ws.Cells.Item(1, 14).Value = "researchArea"
$ws.Cells.Item(1, 15).Value = "leaveExtension"
$ws.Cells.Item(1, 16).Value = "intendedDegree"
$ws.Cells.Item(1, 17).Value = "hoursCompleted"
$ws.Cells.Item(1, 18).Value = "citizenship"
$ws.Cells.Item(1, 19).Value = "fundingEligibility"
$ws.Cells.Item(1, 20).Value = "semestersOnLeave"
$ws.Cells.Item(1, 21).Value = "backgroundApproved"
$ws.Cells.Item(1, 22).Value = "mastersAwarded"
$ws.Cells.Item(1, 23).Value = "prpPassed"
$ws.Cells.Item(1, 24).Value = "technicalWritingApproved"
$ws.Cells.Item(1, 25).Value = "proceedToPhdFormSubmitted"
$ws.Cells.Item(1, 26).Value = "backgroundPrepWorksheetApproved"
$ws.Cells.Item(1, 27).Value = "programOfStudyApproved"
$ws.Cells.Item(1, 28).Value = "researchPlanningMeeting"
$ws.Cells.Item(1, 29).Value = "programProductRequirement"
$ws.Cells.Item(1, 30).Value = "committeeCompApproved"
$ws.Cells.Item(1, 31).Value = "phdProposalApproved"
$ws.Cells.Item(1, 32).Value = "phdAwarded"
$ws.Cells.Item(1, 33).Value = "oralExamPassed"
$ws.Cells.Item(1, 34).Value = "dissertationDefencePassed"
$ws.Cells.Item(1, 35).Value = "dissertationSubmitted"
$ws.Cells.Item(1, 36).Value = "jobHistory"
$ws.Cells.Item(1, 37).Value = "semesterStarted"
$ws.Cells.Item(1, 38).Value = "advisor"
$ws.Cells.Item(1, 39).Value = "otherAdvisor"
$ws.Cells.Item(1, 40).Value = "researchAdvisor"
$ws.Cells.Item(1, 41).Value = "otherResearchAdvisor"
$ws.Cells.Item(1, 42).Value = "grades"
$ws.Cells.Item(1, 43).Value = "courseHistory"

# ---- Row 2 (hmbodnar) ----
$ws.Cells.Item(2, 1).Value = "hmbodnar"
$ws.Cells.Item(2, 2).Value = "iddd"
$ws.Cells.Item(2, 3).Value = "hannahbodnar.17+student@gmail.com"
$ws.Cells.Item(2, 4).Value = "hannah"
$ws.Cells.Item(2, 5).Value = "bodnar"
$ws.Cells.Item(2, 6).Value = "None"
$ws.Cells.Item(2, 7).Value = 730171699
$ws.Cells.Item(2, 8).Value = "Graduated"
$ws.Cells.Item(2, 10).Value = "OTHER"
$ws.Cells.Item(2, 11).Value = "OTHER"
$ws.Cells.Item(2, 12).Value = "NO"
$ws.Cells.Item(2, 16).Value = "MASTERS"
$ws.Cells.Item(2, 17).Value = 30
$ws.Cells.Item(2, 19).Value = "PROBATION"
$ws.Cells.Item(2, 38).Value = "person, person"
$ws.Cells.Item(2, 40).Value = "test, test"

# ---- Row 3 (fakeonyen) ----
$ws.Cells.Item(3, 1).Value = "fakeonyen"
$ws.Cells.Item(3, 2).Value = "fakecsid"
$ws.Cells.Item(3, 3).Value = "fakeEmail@fake.com"
$ws.Cells.Item(3, 4).Value = "fake"
$ws.Cells.Item(3, 5).Value = "fake"
$ws.Cells.Item(3, 6).Value = "she, her"
$ws.Cells.Item(3, 7).Value = 949949949
$ws.Cells.Item(3, 8).Value = "Graduated"
$ws.Cells.Item(3, 9).Value = "fake"
$ws.Cells.Item(3, 10).Value = "FEMALE"
$ws.Cells.Item(3, 11).Value = "OTHER"
$ws.Cells.Item(3, 12).Value = "YES"
$ws.Cells.Item(3, 13).Value = "help"
$ws.Cells.Item(3, 14).Value = "Systems"
$ws.Cells.Item(3, 15).Value = "NO"
$ws.Cells.Item(3, 16).Value = "MASTERS"
$ws.Cells.Item(3, 17).Value = 20
$ws.Cells.Item(3, 18).Value = $true
$ws.Cells.Item(3, 19).Value = "GUARANTEED"
$ws.Cells.Item(3, 20).Value = 99
$ws.Cells.Item(3, 21).NumberFormat = "@"
$ws.Cells.Item(3, 21).Value = "2019-09-18"
$ws.Cells.Item(3, 22).NumberFormat = "@"
$ws.Cells.Item(3, 22).Value = "2019-02-01"
$ws.Cells.Item(3, 23).NumberFormat = "@"
$ws.Cells.Item(3, 23).Value = "2019-09-18"
$ws.Cells.Item(3, 24).NumberFormat = "@"
$ws.Cells.Item(3, 24).Value = "2019-09-18"
$ws.Cells.Item(3, 25).Value = ""
$ws.Cells.Item(3, 26).NumberFormat = "@"
$ws.Cells.Item(3, 26).Value = "2019-09-18"
$ws.Cells.Item(3, 27).NumberFormat = "@"
$ws.Cells.Item(3, 27).Value = "2019-09-18"
$ws.Cells.Item(3, 28).NumberFormat = "@"
$ws.Cells.Item(3, 28).Value = "2019-09-18"
$ws.Cells.Item(3, 29).NumberFormat = "@"
$ws.Cells.Item(3, 29).Value = "2019-09-18"
$ws.Cells.Item(3, 30).NumberFormat = "@"
$ws.Cells.Item(3, 30).Value = "2019-09-18"
$ws.Cells.Item(3, 31).NumberFormat = "@"
$ws.Cells.Item(3, 31).Value = "2019-09-18"
$ws.Cells.Item(3, 32).NumberFormat = "@"
$ws.Cells.Item(3, 32).Value = "2019-02-01"
$ws.Cells.Item(3, 33).NumberFormat = "@"
$ws.Cells.Item(3, 33).Value = "2019-09-18"
$ws.Cells.Item(3, 34).NumberFormat = "@"
$ws.Cells.Item(3, 34).Value = "2019-09-18"
$ws.Cells.Item(3, 35).NumberFormat = "@"
$ws.Cells.Item(3, 35).Value = "2019-09-18"
$ws.Cells.Item(3, 37).Value = "FA 2019"
$ws.Cells.Item(3, 38).Value = "test, test"
$ws.Cells.Item(3, 40).Value = "test, test"

# ---- Row 4 (new) ----
$ws.Cells.Item(4, 1).Value = "new"
$ws.Cells.Item(4, 2).Value = "newnew"
$ws.Cells.Item(4, 3).Value = "new@gmail.com"
$ws.Cells.Item(4, 4).Value = "new"
$ws.Cells.Item(4, 5).Value = "student"
$ws.Cells.Item(4, 6).Value = "None"
$ws.Cells.Item(4, 7).Value = 111111111
$ws.Cells.Item(4, 8).Value = "Active"
$ws.Cells.Item(4, 10).Value = "OTHER"
$ws.Cells.Item(4, 11).Value = "OTHER"
$ws.Cells.Item(4, 12).Value = "NO"
$ws.Cells.Item(4, 16).Value = "MASTERS"
$ws.Cells.Item(4, 18).Value = $false
$ws.Cells.Item(4, 19).Value = "NOT GUARANTEED"
$ws.Cells.Item(4, 38).Value = "test, test"

# ---- Row 5 (fake/upload) ----
$ws.Cells.Item(5, 1).Value = "fake"
$ws.Cells.Item(5, 2).Value = "fake"
$ws.Cells.Item(5, 3).Value = "fake@gmail.com"
$ws.Cells.Item(5, 4).Value = "test"
$ws.Cells.Item(5, 5).Value = "upload"
$ws.Cells.Item(5, 6).Value = "she, her"
$ws.Cells.Item(5, 7).Value = 2828282
$ws.Cells.Item(5, 8).Value = "Graduated"
$ws.Cells.Item(5, 9).Value = "fake"
$ws.Cells.Item(5, 10).Value = "FEMALE"
$ws.Cells.Item(5, 11).Value = "OTHER"
$ws.Cells.Item(5, 12).Value = "YES"
$ws.Cells.Item(5, 13).Value = "help"
$ws.Cells.Item(5, 14).Value = "Systems"
$ws.Cells.Item(5, 15).Value = "NO"
$ws.Cells.Item(5, 16).Value = "MASTERS"
$ws.Cells.Item(5, 17).Value = 20
$ws.Cells.Item(5, 18).Value = $true
$ws.Cells.Item(5, 19).Value = "GUARANTEED"
$ws.Cells.Item(5, 20).Value = 99
$ws.Cells.Item(5, 21).NumberFormat = "@"
$ws.Cells.Item(5, 21).Value = "2019-09-18"
$ws.Cells.Item(5, 22).NumberFormat = "@"
$ws.Cells.Item(5, 22).Value = "2019-02-01"
$ws.Cells.Item(5, 23).NumberFormat = "@"
$ws.Cells.Item(5, 23).Value = "2019-09-18"
$ws.Cells.Item(5, 24).NumberFormat = "@"
$ws.Cells.Item(5, 24).Value = "2019-09-18"
$ws.Cells.Item(5, 25).Value = ""
$ws.Cells.Item(5, 26).NumberFormat = "@"
$ws.Cells.Item(5, 26).Value = "2019-09-18"
$ws.Cells.Item(5, 27).NumberFormat = "@"
$ws.Cells.Item(5, 27).Value = "2019-09-18"
$ws.Cells.Item(5, 28).NumberFormat = "@"
$ws.Cells.Item(5, 28).Value = "2019-09-18"
$ws.Cells.Item(5, 29).NumberFormat = "@"
$ws.Cells.Item(5, 29).Value = "2019-09-18"
$ws.Cells.Item(5, 30).NumberFormat = "@"
$ws.Cells.Item(5, 30).Value = "2019-09-18"
$ws.Cells.Item(5, 31).NumberFormat = "@"
$ws.Cells.Item(5, 31).Value = "2019-09-18"
$ws.Cells.Item(5, 32).NumberFormat = "@"
$ws.Cells.Item(5, 32).Value = "2019-02-01"
$ws.Cells.Item(5, 33).NumberFormat = "@"
$ws.Cells.Item(5, 33).Value = "2019-09-18"
$ws.Cells.Item(5, 34).NumberFormat = "@"
$ws.Cells.Item(5, 34).Value = "2019-09-18"
$ws.Cells.Item(5, 35).NumberFormat = "@"
$ws.Cells.Item(5, 35).Value = "2019-09-18"
$ws.Cells.Item(5, 37).Value = "FA 2019"
$ws.Cells.Item(5, 39).Value = "test, test"
$ws.Cells.Item(5, 40).Value = "Pozefsky, Diane"
